# "multithreading + crm api connected"
# - Update the listing description text for the Hollywood, FL villa (C2)
# - Give that cell a distinct font + wrap-text formatting so the longer
#   description still reads well, which also bumps rows 2 & 3 a bit taller
# - Nudge column widths back to their natural (auto) values
# - Leave the selection on C10, like the author did before saving

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New listing description text (same shared-string slot, new wording)
$ws.Range("C2").Value = "Top Villa w/Game room HotTub & HTD Pool Near Beach"

# 2. Re-stamp the font (forces a fresh font record) and turn wrapping on
#    for that cell so the longer text is readable in the column.
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").WrapText = $true

# 3. Rows holding data grew slightly taller once the new font/wrap kicked in
$ws.Rows.Item(2).RowHeight = 16
$ws.Rows.Item(3).RowHeight = 16

# 4. Re-apply column widths (values unchanged) so the stored widths get
#    recomputed/normalised the way Excel does on save
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# 5. Leave the cursor one row further down, like in the saved file
$ws.Range("C10").Select() | Out-Null
